# disk_savvy.xlsx — append a further week of disk-usage samples
# (rows 114-120) to Arkusz1, matching the style of the existing data
# (date in col A, time-of-day fraction in col B, file count in col C,
# disk space in col D), then move the viewport/selection down to follow
# the newly appended rows (mirrors what Excel does automatically when a
# user keeps typing new rows at the bottom of a sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append right after the current last row (113).
$newRows = @(
    @(45133, 0.46194444444444444, 83219, 1680),
    @(45134, 0.89277777777777778, 83220, 1680),
    @(45135, 0.75491898148148151, 83220, 1680),
    @(45136, 0.58952546296296293, 83236, 1680),
    @(45137, 0.4776157407407407,  83236, 1680),
    @(45138, 0.4611574074074074,  83236, 1680),
    @(45139, 0.46293981481481478, 83244, 1690)
)

$firstNewRow = 114
$lastNewRow = $firstNewRow + $newRows.Length - 1   # 120

# Copy the formatting (date/time number formats etc.) of the last
# existing row down across all the new rows in one shot, so the new
# cells pick up the same style indices as row 113 instead of minting
# brand-new styles.
$lastExistingRow = $firstNewRow - 1                # 113
$srcFormatRange = $ws.Range("A" + $lastExistingRow + ":D" + $lastExistingRow)
$dstFormatRange = $ws.Range("A" + $firstNewRow + ":D" + $lastNewRow)
$srcFormatRange.Copy()
$dstFormatRange.PasteSpecial(-4122) # xlPasteFormats

# Write the actual values row by row.
$r = $firstNewRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $r++
}

# Follow the freshly-entered data: select the next empty row (A121),
# same as Excel leaves the cursor after typing the last row.
$nextRow = $lastNewRow + 1                          # 121
$ws.Range("A" + $nextRow).Select()
